# Edits QE_holdings.xlsx per commit "Add files via upload":
#  - Bumps the "as of" date in the confidential disclosure footnote
#    from 2021-03-30 to 2021-03-31 (cell A44).
#  - Refreshes the Weight (D) and Percent Change (E) values for every
#    holding row (2-41) with the new 2021-03-31 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect so the locked data cells can be updated.
$ws.Unprotect()

# Update the confidential "as of" date disclosure text.
$ws.Range("A44").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-31 for illustrative purposes only and are subject to change."
# Re-run autofit on the row so the wrapped two-line text doesn't leave behind an
# explicit/custom row height (keeps row 44 identical in shape to the source file).
$ws.Rows.Item(44).AutoFit()

# Refresh Weight (D) and Percent Change (E) for each holding row.
$ws.Range("D2").Value = 0.07634629361622003
$ws.Range("E2").Value = 0.01876563803169318
$ws.Range("D3").Value = 0.06678519308222956
$ws.Range("E3").Value = 0.01690748328660785
$ws.Range("D4").Value = 0.05480080497295482
$ws.Range("E4").Value = -0.01053487741585657
$ws.Range("D5").Value = 0.04632035970902357
$ws.Range("E5").Value = 0.01269601248981278
$ws.Range("D6").Value = 0.04449849742222481
$ws.Range("E6").Value = -0.01456499223200414
$ws.Range("D7").Value = 0.03855064843864416
$ws.Range("E7").Value = -0.0008022652194432656
$ws.Range("D8").Value = 0.04002669497474075
$ws.Range("E8").Value = -0.003999757590449016
$ws.Range("D9").Value = 0.03498454669028098
$ws.Range("E9").Value = 0.0006630322675704736
$ws.Range("D10").Value = 0.03102593421041062
$ws.Range("E10").Value = 0.007842773165499528
$ws.Range("D11").Value = 0.02812766165180766
$ws.Range("E11").Value = -0.005443863526114368
$ws.Range("D12").Value = 0.03192026437849133
$ws.Range("E12").Value = -0.007694280584765401
$ws.Range("D13").Value = 0.03056398743601802
$ws.Range("E13").Value = 0.02267361111111121
$ws.Range("D14").Value = 0.02715282812396244
$ws.Range("E14").Value = -0.006002233389168055
$ws.Range("D15").Value = 0.03080276858786191
$ws.Range("E15").Value = -0.000147655961609261
$ws.Range("D16").Value = 0.02831717050247741
$ws.Range("E16").Value = -0.003988649748367101
$ws.Range("D17").Value = 0.02772772217335421
$ws.Range("E17").Value = 0.001410668591299791
$ws.Range("D18").Value = 0.02320316046184368
$ws.Range("E18").Value = 0.003606711619883907
$ws.Range("D19").Value = 0.01927289850603351
$ws.Range("E19").Value = 0.05083225826751825
$ws.Range("D20").Value = 0.02250546467721791
$ws.Range("E20").Value = -0.01582393597671872
$ws.Range("D21").Value = 0.02148651348893678
$ws.Range("E21").Value = -0.01517022402540136
$ws.Range("D22").Value = 0.02220664712148184
$ws.Range("E22").Value = -0.007509813961426914
$ws.Range("D23").Value = 0.02095058245924271
$ws.Range("E23").Value = -0.008278457196613243
$ws.Range("D24").Value = 0.02049909657340709
$ws.Range("E24").Value = -0.01496908558411991
$ws.Range("D25").Value = 0.01815494789416249
$ws.Range("E25").Value = 0.0007515657620043026
$ws.Range("D26").Value = 0.01814054522151159
$ws.Range("E26").Value = -0.0004512974802558523
$ws.Range("D27").Value = 0.01962174639834639
$ws.Range("E27").Value = -0.001158972377825029
$ws.Range("D28").Value = 0.01766571184527345
$ws.Range("E28").Value = -0.01076182384593594
$ws.Range("D29").Value = 0.01866828946885672
$ws.Range("E29").Value = 0.001689189189189255
$ws.Range("D30").Value = 0.01806595453788797
$ws.Range("E30").Value = 0.003323179174743895
$ws.Range("D31").Value = 0.01818299520406162
$ws.Range("E31").Value = -0.005386250885896393
$ws.Range("D32").Value = 0.01619012013041848
$ws.Range("E32").Value = 0.01339076692574204
$ws.Range("D33").Value = 0.01728805860765877
$ws.Range("E33").Value = -0.007646976287357998
$ws.Range("D34").Value = 0.007805793755546925
$ws.Range("E34").Value = 0.03701905335327349
$ws.Range("D35").Value = 0.007783355907627625
$ws.Range("E35").Value = 0.01610861138705455
$ws.Range("D36").Value = 0.007172227765987801
$ws.Range("E36").Value = 0.02663397311237015
$ws.Range("D37").Value = 0.006367497382503755
$ws.Range("E37").Value = 0.008904761904761971
$ws.Range("D38").Value = 0.006915708585721221
$ws.Range("E38").Value = 0.01644159943879342
$ws.Range("D39").Value = 0.007056703170619519
$ws.Range("E39").Value = 0.02129076612383463
$ws.Range("D40").Value = 0.006844604864949929
$ws.Range("E40").Value = 0.0108755842027155
$ws.Range("D41").Value = 1
$ws.Range("E41").Value = 0.00265843015719569
